$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 88, shifting existing rows 88:197 down to 89:198.
$ws.Rows.Item(88).Insert()

# Populate the newly inserted row 88 with its data.
$ws.Range("A88").Value = 1
$ws.Range("B88").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C88").Value = "Arica y Parinacota"
$ws.Range("D88").Value = 44482
$ws.Range("D88").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E88").Value = 15
$ws.Range("F88").Value = 100114013
$ws.Range("G88").Value = "Zanahoria"
$ws.Range("H88").Value = "Sin especificar"
$ws.Range("I88").Value = "Primera"
$ws.Range("J88").Value = 80
$ws.Range("K88").Value = 8000
$ws.Range("L88").Value = 9000
$ws.Range("M88").Value = 8500
$ws.Range("N88").Value = "$/saco 25 kilos"
$ws.Range("O88").Value = "Valle de Camiña"
$ws.Range("P88").Value = 340
$ws.Range("Q88").Value = 25
$ws.Range("R88").Value = "Hortaliza"
